$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N (14th column) to make room
# for the "Variable Instalments" related data - this shifts the old
# N/O/P columns (Due / In Advance / Over Due) one column to the right.
$ws.Columns("N:N").Insert()

# Match the new column's width to its neighbour (column M), mirroring
# what Excel does when a column is inserted next to formatted data.
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Make "Repayment Schedule" the active sheet / tab.
$ws.Activate()

# Restore the previously-selected cell on this sheet.
$ws.Range("L15").Select() | Out-Null
